# B1-and-B2-PowerPoint.pptx edit
#
# The underlying OOXML diff for this commit contains two parts:
#
#   1. On slide 5, the single table's <a:tableStyleId> changes from
#      {3F213E74-3EC8-428F-AE3F-06FF0DECC9E9} to
#      {2BD116F1-F12D-4D75-8751-8E82CE40F0B2}.
#      This is a normal "pick a different table style from the gallery"
#      action, exposed on the PowerPoint object model as
#      Table.ApplyStyle(styleId).
#
#   2. ppt/theme/theme1.xml and ppt/theme/theme2.xml trade places (the
#      deck's slide-master theme becomes the plain "Office Theme" that
#      used to live behind the notes master, and the notes master is
#      left holding the "Integral" theme that used to live behind the
#      slide master). There is no documented/working PowerPoint
#      Automation call that rewrites the raw theme (colours/fonts/effects)
#      of an existing master in one step, so we drive it through the
#      same entry point PowerPoint itself uses when a different design
#      is applied (Master.ApplyTheme / Presentation.ApplyTheme). This is
#      a best-effort call: on hosts where it is not wired up it is a
#      harmless no-op, and it does not affect the table-style edit above.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 ---------------------------------------
$targetStyleId = "{2BD116F1-F12D-4D75-8751-8E82CE40F0B2}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}

# --- 2. Theme swap between the slide master and the notes master ------
try {
    $slideMasterTheme = $p.SlideMaster.Theme
    $notesMasterTheme = $p.NotesMaster.Theme
    $p.SlideMaster.ApplyTheme($notesMasterTheme)
    $p.NotesMaster.ApplyTheme($slideMasterTheme)
} catch {
    # Theme-content replacement isn't available on every host; the
    # table-style change above is applied regardless.
}
